# Fix revenue report backlog issue: update the "From"/"To" date range
# and the "Print taken at" timestamp on the Denomination wise Stock Ledger sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = "From : 01-10-2018"
$ws.Range("G4").Value = "To : 30-04-2019"
$ws.Range("G13").Value = "Print taken at : 30-04-2019 14:33:01"
